# Taller_serie_Taylor.xlsx -- "Add files via upload" edit
# Adds a second Taylor-series table (columns G:K) mirroring the existing
# one (columns A:E), restyles both tables with a blue banded look, fixes
# the A1 label to use a comma decimal separator, widens a couple of
# columns and tweaks the view (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Color constants (BGR-packed ints, as COM ColorIndex/Color expects)
# Derived from theme accent5 (#5B9BD5) tinted like the uploaded file:
#   light band  -> tint -0.25 -> #4474A0
#   dark band   -> tint -0.50 -> #2E4E6B
# ------------------------------------------------------------------
$white = 16777215
$lightBlue = 10515524
$darkBlue  = 7032366

# ------------------------------------------------------------------
# 1) Fix the A1 label (decimal comma, like the rest of the sheet)
# ------------------------------------------------------------------
$ws.Range("A1").Value = "f(0,6)"

# ------------------------------------------------------------------
# 2) Add the second table (columns G:K), mirroring A:E
# ------------------------------------------------------------------
$ws.Range("G1").Value = "f(0,45)"
$ws.Range("H1").Value = "(1,6*EXP(1)^(0,45))-(4,2(0,45))+2,75"
$ws.Range("K1").Formula = "=(1.6*EXP(1)^(0.45))-(4.2*(0.45))+2.75"

$ws.Range("G2").Value = "orden [0]"
$ws.Range("H2").Value = "(1,6*EXP(1)^(0,4))-(4,2(0,4))+2,75"
$ws.Range("K2").Formula = "=(1.6*EXP(1)^(0.4))-(4.2*(0.4))+2.75"

$ws.Range("G3").Value = "orden [1]"
$ws.Range("H3").Value = "3,4569195+(((8*EXP(1)^(0,4))-21)/5)*(0,05)"
$ws.Range("K3").Formula = "=3.4569195+(((8*EXP(1)^(0.4))-21)/5)*(0.05)"

$ws.Range("G4").Value = "orden [2]"
$ws.Range("H4").Value = "3,3662655+((8*EXP(1)^(0,4))/2)*(0,05)^2"
$ws.Range("K4").Formula = "=3.3662655+((8*EXP(1)^(0.4))/2)*(0.05)^2"

$ws.Range("G5").Value = "orden [3]"
$ws.Range("H5").Value = "3,3811837+((8*EXP(1)^(0,4))/6)*(0,05)^3"
$ws.Range("K5").Formula = "=3.3811837+((8*EXP(1)^(0.4))/6)*(0.05)^3"

# Merge the label columns of the new table, same as H:J under B:D
$ws.Range("H1:J1").Merge()
$ws.Range("H2:J2").Merge()
$ws.Range("H3:J3").Merge()
$ws.Range("H4:J4").Merge()
$ws.Range("H5:J5").Merge()

# ------------------------------------------------------------------
# 3) Style both tables with the blue banded look (header row lighter,
#    body rows darker; white bold-less font, thin white borders,
#    centered both ways)
# ------------------------------------------------------------------

# Header row (row 1) - light band, for first column + merged block
$headerFirstCol = $ws.Range("A1,G1")
$headerFirstCol.Font.Color = $white
$headerFirstCol.Interior.Color = $lightBlue
$headerFirstCol.Borders.LineStyle = 1
$headerFirstCol.Borders.Color = $white
$headerFirstCol.HorizontalAlignment = -4108
$headerFirstCol.VerticalAlignment = -4108

$headerLabel = $ws.Range("B1:D1,H1:J1")
$headerLabel.Font.Color = $white
$headerLabel.Interior.Color = $lightBlue
$headerLabel.Borders.LineStyle = 1
$headerLabel.Borders.Color = $white
$headerLabel.HorizontalAlignment = -4108
$headerLabel.VerticalAlignment = -4108

$headerResult = $ws.Range("E1,K1")
$headerResult.Font.Color = $white
$headerResult.Interior.Color = $lightBlue
$headerResult.Borders.LineStyle = 1
$headerResult.Borders.Color = $white
$headerResult.HorizontalAlignment = -4108
$headerResult.VerticalAlignment = -4108

# Body rows (2-5) - dark band, for first column + merged block + result,
# except the last row's result cell which (like the header) uses the
# lighter band in the source file.
$bodyFirstCol = $ws.Range("A2:A5,G2:G5")
$bodyFirstCol.Font.Color = $white
$bodyFirstCol.Interior.Color = $darkBlue
$bodyFirstCol.Borders.LineStyle = 1
$bodyFirstCol.Borders.Color = $white
$bodyFirstCol.HorizontalAlignment = -4108
$bodyFirstCol.VerticalAlignment = -4108

$bodyLabel = $ws.Range("B2:D5,H2:J5")
$bodyLabel.Font.Color = $white
$bodyLabel.Interior.Color = $darkBlue
$bodyLabel.Borders.LineStyle = 1
$bodyLabel.Borders.Color = $white
$bodyLabel.HorizontalAlignment = -4108
$bodyLabel.VerticalAlignment = -4108

$bodyResultDark = $ws.Range("E2:E4,K2:K4")
$bodyResultDark.Font.Color = $white
$bodyResultDark.Interior.Color = $darkBlue
$bodyResultDark.Borders.LineStyle = 1
$bodyResultDark.Borders.Color = $white
$bodyResultDark.HorizontalAlignment = -4108
$bodyResultDark.VerticalAlignment = -4108

$bodyResultLight = $ws.Range("E5,K5")
$bodyResultLight.Font.Color = $white
$bodyResultLight.Interior.Color = $lightBlue
$bodyResultLight.Borders.LineStyle = 1
$bodyResultLight.Borders.Color = $white
$bodyResultLight.HorizontalAlignment = -4108
$bodyResultLight.VerticalAlignment = -4108

# ------------------------------------------------------------------
# 4) Column widths
# ------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 15.75
$ws.Columns.Item(10).ColumnWidth = 22.25

# ------------------------------------------------------------------
# 5) View tweaks: zoom in to 140% and move the selection
# ------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 140
$ws.Range("G7").Select()
